# Updates cryptos list values (prices + 1h volume change) per commit
# "Updated cryptos list on Mon Jul  3 04:17:59 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.747.43'
$ws.Range('E2').Value = '  +0.67%  '
$ws.Range('D3').Value = '1.946.86'
$ws.Range('E3').Value = '  +1.56%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4824'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.60%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2957'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.91%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06842'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '112.72'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.51'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.49%  '
$ws.Range('D12').Value = '1.948.69'
$ws.Range('E12').Value = '  +1.68%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.562'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07657'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6933'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '297.52'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +8.91%  '
$ws.Range('D17').Value = '30.758.15'
$ws.Range('E17').Value = '  +0.71%  '
$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.43'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.78%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.733'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.33%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007720'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.61%  '
$ws.Range('D21').Value = '2.201.13'
$ws.Range('E21').Value = '  +1.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9996'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.578'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.792'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.81'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.44'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.42%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.183'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1093'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.442'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.774'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +18.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.457'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +8.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05075'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7800'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +7.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.165'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02072'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.735'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.705'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.042'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '111.15'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4473'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8765'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.978'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.79%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '71.39'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.36%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9997'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.421'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.469'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '48.98'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1259'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '35.71'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.2566'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.93%  '
